# "Corporate APNs DB GUI"
#
# Adds two new header columns ("Destination" and "Corporate APNs") to the
# right of the existing "MTX Name" / "MTX Number" table on Sheet1, and
# updates the sheet selection to cover the whole grid (as if the user had
# pressed Ctrl+A / "Select All" after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1. Setting .Value adds the strings to the shared
# string table, extends <dimension>, and grows every row's "spans"
# attribute automatically - the header cells pick up the same style (bold,
# s="1") that the row itself already carries, matching A1/B1.
$ws.Range("C1").Value = "Destination"
$ws.Range("D1").Value = "Corporate APNs"

# Size the two new columns to (approximately) the same "best fit" pixel
# widths recorded in the target workbook (11.36328125 and 15.54296875
# characters). The engine only supports whole-pixel column widths, so we
# pick the ColumnWidth input that rounds to the closest achievable pixel
# width (68 px and 93 px respectively).
$ws.Columns.Item(3).ColumnWidth = 10.5
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666

# Select the entire sheet (A1:XFD1048576), matching the saved sheet view.
$ws.Cells.Select()
